# Update the Fin Buff Calc inputs with the latest figures from 502 Part C / Part L
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Gross Expenditures From 502 Part C
$ws.Range("D3").Value = 740798.39

# Total Labor Cost From 502 Part L
$ws.Range("D5").Value = 320631.03000000003

# Leave the selection where the user last worked in the sheet
$ws.Range("G7:I16").Select() | Out-Null
